# August additions to job search
# 1) Add a new row to the "Big Companies" sheet for the Netflix application.
# 2) Add a new "Interviews" sheet tracking companies that followed up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Big Companies: append the Netflix application as row 36
# ---------------------------------------------------------------------------
$bigCompanies = $wb.Worksheets.Item("Big Companies")

$bigCompanies.Range("A36").Value = "Software Engineer, Personalization Application Core"
$bigCompanies.Range("B36").Value = "Netflix"
$bigCompanies.Range("C36").Value = "Submitted 8/4/22"
$bigCompanies.Range("D36").Value = "https://jobs.netflix.com/jobs/219156808"

[void]$bigCompanies.Range("D36").Select()

# ---------------------------------------------------------------------------
# 2. New "Interviews" sheet, placed after "Big Companies"
# ---------------------------------------------------------------------------
$interviews = $wb.Worksheets.Add($null, $bigCompanies)
$interviews.Name = "Interviews"

# Column widths (values chosen so the stored OOXML <col> widths land on
# 23.42578125 / 29 / 19 / 31.42578125 after the host's pixel-quantized
# ColumnWidth -> stored-width conversion)
$interviews.Columns.Item(1).ColumnWidth = 22.666666666666668
$interviews.Columns.Item(2).ColumnWidth = 28.166666666666668
$interviews.Columns.Item(4).ColumnWidth = 18.166666666666668
$interviews.Columns.Item(5).ColumnWidth = 30.666666666666668

# Cell contents (order chosen to match the original authoring sequence)
$interviews.Range("A3").Value = "Enterprise Fleet"
$interviews.Range("A1").Value = "Company"
$interviews.Range("B1").Value = "Job Tech Stack"
$interviews.Range("C1").Value = "Salary"
$interviews.Range("C2").Value = "71k"
$interviews.Range("C3").Value = "125k"
$interviews.Range("D1").Value = "Interview Time"
$interviews.Range("B3").Value = "Java,Angular,hosted on Oracle"
$interviews.Range("E1").Value = "How I found"
$interviews.Range("E2").Value = "Recruiter Linked Reached Out to Me"
$interviews.Range("E3").Value = "Recruiter Linked Reached Out to Me"
$interviews.Range("A2").Value = "Oracle Cerner"

# Header row formatting - bold
$interviews.Range("A1:E1").Font.Bold = $true

# Interview date (8/9/2022) in D2, formatted as a short date (built-in date format)
$interviews.Range("D2").NumberFormat = "mm-dd-yy"
$interviews.Range("D2").Value = Get-Date -Year 2022 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0

[void]$interviews.Range("E3").Select()
[void]$interviews.Activate()

$wb.Save() | Out-Null
